$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values of the columns that rotate between rows 2, 3 and 4
# (D = Fecha, I = Calidad, J = Volumen, K = Precio minimo, L = Precio maximo,
#  M = Precio promedio ponderado, P = Precio $/Kg)
$row2 = @{
    D = $ws.Range("D2").Value2
    I = $ws.Range("I2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    P = $ws.Range("P2").Value2
}
$row3 = @{
    D = $ws.Range("D3").Value2
    I = $ws.Range("I3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    P = $ws.Range("P3").Value2
}
$row4 = @{
    D = $ws.Range("D4").Value2
    I = $ws.Range("I4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}

# Row 2 takes on the old Row 3 values
$ws.Range("D2").Value2 = $row3.D
$ws.Range("I2").Value2 = $row3.I
$ws.Range("J2").Value2 = $row3.J
$ws.Range("K2").Value2 = $row3.K
$ws.Range("L2").Value2 = $row3.L
$ws.Range("M2").Value2 = $row3.M
$ws.Range("P2").Value2 = $row3.P

# Row 3 takes on the old Row 4 values
$ws.Range("D3").Value2 = $row4.D
$ws.Range("I3").Value2 = $row4.I
$ws.Range("J3").Value2 = $row4.J
$ws.Range("K3").Value2 = $row4.K
$ws.Range("L3").Value2 = $row4.L
$ws.Range("M3").Value2 = $row4.M
$ws.Range("P3").Value2 = $row4.P

# Row 4 takes on the old Row 2 values
$ws.Range("D4").Value2 = $row2.D
$ws.Range("I4").Value2 = $row2.I
$ws.Range("J4").Value2 = $row2.J
$ws.Range("K4").Value2 = $row2.K
$ws.Range("L4").Value2 = $row2.L
$ws.Range("M4").Value2 = $row2.M
$ws.Range("P4").Value2 = $row2.P
